$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("A2").Value = "05 Nov 2025, 01:02 PM"

$ws7 = $wb.Worksheets.Item("Industry Analysis")
$ws7.Range("F2").Value = 21.3
$ws7.Range("F3").Value = -4.3927
$ws7.Range("F4").Value = 35.9445
$ws7.Range("F5").Value = -51.0482
$ws7.Range("F6").Value = 57.2275
$ws7.Range("F7").Value = -9.640700000000001
$ws7.Range("F8").Value = -6.1449
$ws7.Range("F9").Value = 36.9733
$ws7.Range("F10").Value = -4.7026
$ws7.Range("F11").Value = 46.5317
$ws7.Range("F12").Value = -2.102
$ws7.Range("F13").Value = 17.4681
$ws7.Range("F14").Value = -33.0245
$ws7.Range("F15").Value = 1.0205
$ws7.Range("F16").Value = 2.0426
$ws7.Range("F17").Value = -16.2411
$ws7.Range("F18").Value = 7.4627
$ws7.Range("F19").Value = -25.798
$ws7.Range("F20").Value = 47.7485
$ws7.Range("F21").Value = 19.5587
$ws7.Range("F22").Value = 76.5603
$ws7.Range("F23").Value = -54.2675
$ws7.Range("F24").Value = -0.8811
$ws7.Range("F25").Value = 4.8518
$ws7.Range("F26").Value = 3.6831
$ws7.Range("F27").Value = -34.0874
$ws7.Range("F28").Value = -11.9893
$ws7.Range("F29").Value = -12.994
$ws7.Range("F30").Value = 25.5415
$ws7.Range("F31").Value = 56.5088
$ws7.Range("F32").Value = 2.0908
$ws7.Range("F33").Value = -4.7193
$ws7.Range("F34").Value = 22.8807
$ws7.Range("F35").Value = 5.3359
$ws7.Range("F36").Value = -5.1995
$ws7.Range("F37").Value = -5.6238
$ws7.Range("F38").Value = -22.595
$ws7.Range("F39").Value = 10.8405
$ws7.Range("F40").Value = -7.5963
$ws7.Range("F41").Value = -4.552
$ws7.Range("F42").Value = 22.3098
$ws7.Range("F43").Value = 14.0694
$ws7.Range("F44").Value = -9.6066
$ws7.Range("F45").Value = 27.639
$ws7.Range("F46").Value = -6.3484
$ws7.Range("F47").Value = -40.5302
$ws7.Range("F48").Value = -29.7988
$ws7.Range("F49").Value = -24.0791
$ws7.Range("F50").Value = -49.1803
$ws7.Range("F51").Value = -51.6023
$ws7.Range("F52").Value = -34.4756
$ws7.Range("F53").Value = -11.5478
$ws7.Range("F54").Value = -2.3796
$ws7.Range("F55").Value = -15.4382
$ws7.Range("F56").Value = -27.6987
$ws7.Range("F57").Value = -27.1559
$ws7.Range("F58").Value = -2.1585
$ws7.Range("F59").Value = -23.0964
$ws7.Range("F60").Value = -13.3217
$ws7.Range("F61").Value = -8.1496
$ws7.Range("F62").Value = -16.0695
$ws7.Range("F63").Value = -12.5465
$ws7.Range("F64").Value = 47.7264
$ws7.Range("F65").Value = -42.4232
$ws7.Range("F66").Value = 11.3291
$ws7.Range("F67").Value = 14.3746
$ws7.Range("F68").Value = 32.6702
$ws7.Range("F69").Value = -17.0097
$ws7.Range("F70").Value = -13.5162
$ws7.Range("F71").Value = 11.4259
$ws7.Range("F72").Value = 2.6754
$ws7.Range("F73").Value = -11.1574
$ws7.Range("F74").Value = -13.2502
$ws7.Range("F75").Value = 24.7078
$ws7.Range("F76").Value = 53.3554
